$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

$ws.Cells.Item(2, 4).Value = '28.052.83'
$ws.Cells.Item(2, 5).Value = '  +0.23%  '
$ws.Cells.Item(3, 4).Value = '1.847.54'
$ws.Cells.Item(3, 5).Value = '  -0.77%  '
Set-TextValue $ws.Cells.Item(4, 4) '1.005'
$ws.Cells.Item(4, 5).Value = '  +0.31%  '
Set-TextValue $ws.Cells.Item(5, 4) '331.15'
$ws.Cells.Item(5, 5).Value = '  -1.39%  '
Set-TextValue $ws.Cells.Item(6, 4) '1.007'
$ws.Cells.Item(6, 5).Value = '  +0.58%  '
Set-TextValue $ws.Cells.Item(7, 4) '0.4531'
$ws.Cells.Item(7, 5).Value = '  -3.51%  '
Set-TextValue $ws.Cells.Item(8, 4) '0.3881'
$ws.Cells.Item(8, 5).Value = '  -0.49%  '
Set-TextValue $ws.Cells.Item(9, 4) '47.56'
$ws.Cells.Item(9, 5).Value = '  +1.55%  '
Set-TextValue $ws.Cells.Item(10, 4) '0.07725'
$ws.Cells.Item(10, 5).Value = '  -3.02%  '
Set-TextValue $ws.Cells.Item(11, 4) '0.9739'
$ws.Cells.Item(11, 5).Value = '  -0.54%  '
Set-TextValue $ws.Cells.Item(12, 4) '21.20'
$ws.Cells.Item(12, 5).Value = '  -1.45%  '
$ws.Cells.Item(13, 4).Value = '1.846.74'
$ws.Cells.Item(13, 5).Value = '  +0.47%  '
Set-TextValue $ws.Cells.Item(14, 4) '5.755'
$ws.Cells.Item(14, 5).Value = '  -2.97%  '
Set-TextValue $ws.Cells.Item(15, 4) '6.970'
$ws.Cells.Item(15, 5).Value = '  -3.18%  '
Set-TextValue $ws.Cells.Item(16, 4) '1.004'
$ws.Cells.Item(16, 5).Value = '  +0.10%  '
Set-TextValue $ws.Cells.Item(17, 4) '86.85'
$ws.Cells.Item(17, 5).Value = '  -5.12%  '
Set-TextValue $ws.Cells.Item(18, 4) '0.06539'
$ws.Cells.Item(18, 5).Value = '  -1.27%  '
Set-TextValue $ws.Cells.Item(19, 4) '0.00001015'
$ws.Cells.Item(19, 5).Value = '  -2.25%  '
Set-TextValue $ws.Cells.Item(20, 4) '16.89'
$ws.Cells.Item(20, 5).Value = '  -3.44%  '
Set-TextValue $ws.Cells.Item(21, 4) '1.015'
$ws.Cells.Item(21, 5).Value = '  +1.34%  '
$ws.Cells.Item(22, 4).Value = '28.009.51'
$ws.Cells.Item(22, 5).Value = '  +0.13%  '
Set-TextValue $ws.Cells.Item(23, 4) '5.275'
$ws.Cells.Item(23, 5).Value = '  -2.16%  '
Set-TextValue $ws.Cells.Item(24, 4) '10.55'
$ws.Cells.Item(24, 5).Value = '  -3.55%  '
Set-TextValue $ws.Cells.Item(25, 4) '2.257'
$ws.Cells.Item(25, 5).Value = '  -1.45%  '
$ws.Cells.Item(26, 4).Value = '2.066.68'
$ws.Cells.Item(26, 5).Value = '  -0.04%  '
Set-TextValue $ws.Cells.Item(27, 4) '156.04'
$ws.Cells.Item(27, 5).Value = '  -2.04%  '
Set-TextValue $ws.Cells.Item(28, 4) '19.06'
$ws.Cells.Item(28, 5).Value = '  -2.46%  '
Set-TextValue $ws.Cells.Item(29, 4) '2.025'
$ws.Cells.Item(29, 5).Value = '  -3.48%  '
Set-TextValue $ws.Cells.Item(30, 4) '5.222'
$ws.Cells.Item(30, 5).Value = '  -4.38%  '
Set-TextValue $ws.Cells.Item(31, 4) '116.11'
$ws.Cells.Item(31, 5).Value = '  -2.67%  '
Set-TextValue $ws.Cells.Item(32, 4) '0.09240'
$ws.Cells.Item(32, 5).Value = '  -2.52%  '
Set-TextValue $ws.Cells.Item(33, 4) '0.9295'
$ws.Cells.Item(33, 5).Value = '  -2.94%  '
Set-TextValue $ws.Cells.Item(34, 4) '3.616'
$ws.Cells.Item(34, 5).Value = '  +1.06%  '
Set-TextValue $ws.Cells.Item(35, 4) '1.363'
$ws.Cells.Item(35, 5).Value = '  +1.24%  '
Set-TextValue $ws.Cells.Item(36, 4) '5.143'
$ws.Cells.Item(36, 5).Value = '  -3.03%  '
Set-TextValue $ws.Cells.Item(37, 4) '0.05985'
$ws.Cells.Item(37, 5).Value = '  -1.61%  '
Set-TextValue $ws.Cells.Item(38, 4) '0.02179'
$ws.Cells.Item(38, 5).Value = '  -2.95%  '
Set-TextValue $ws.Cells.Item(39, 4) '8.104'
$ws.Cells.Item(39, 5).Value = '  -2.20%  '
Set-TextValue $ws.Cells.Item(40, 4) '1.154'
$ws.Cells.Item(40, 5).Value = '  -0.80%  '
Set-TextValue $ws.Cells.Item(41, 4) '1.010'
$ws.Cells.Item(41, 5).Value = '  +0.83%  '
Set-TextValue $ws.Cells.Item(42, 4) '0.5635'
$ws.Cells.Item(42, 5).Value = '  -4.69%  '
Set-TextValue $ws.Cells.Item(43, 4) '0.1780'
$ws.Cells.Item(43, 5).Value = '  -4.58%  '
Set-TextValue $ws.Cells.Item(44, 4) '9.875'
$ws.Cells.Item(44, 5).Value = '  -3.12%  '
Set-TextValue $ws.Cells.Item(45, 4) '1.242'
$ws.Cells.Item(45, 5).Value = '  -3.04%  '
Set-TextValue $ws.Cells.Item(46, 4) '2.270'
$ws.Cells.Item(46, 5).Value = '  +23.40%  '
Set-TextValue $ws.Cells.Item(47, 4) '0.07161'
$ws.Cells.Item(47, 5).Value = '  +4.14%  '
$ws.Cells.Item(48, 2).Value = 'Decentraland'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue $ws.Cells.Item(48, 4) '0.5334'
$ws.Cells.Item(48, 5).Value = '  -3.69%  '
$ws.Cells.Item(49, 2).Value = 'EnergySwap'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Cells.Item(49, 4) '11.68'
$ws.Cells.Item(49, 5).Value = '  -3.75%  '
Set-TextValue $ws.Cells.Item(50, 4) '1.864'
$ws.Cells.Item(50, 5).Value = '  -4.29%  '
Set-TextValue $ws.Cells.Item(51, 4) '109.40'
$ws.Cells.Item(51, 5).Value = '  -1.85%  '
